$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "37.353.09"
Set-TextValue 2 5 "  +4.50%  "

Set-TextValue 3 4 "2.046.90"
Set-TextValue 3 5 "  +3.00%  "

Set-TextValue 4 5 "  -0.07%  "

Set-TextValue 5 4 "253.09"
Set-TextValue 5 5 "  +3.11%  "

Set-TextValue 6 4 "0.652"
Set-TextValue 6 5 "  +2.07%  "

Set-TextValue 7 4 "65.61"
Set-TextValue 7 5 "  +10.43%  "

Set-TextValue 8 5 "  -0.01%  "

Set-TextValue 9 4 "0.400"
Set-TextValue 9 5 "  +9.60%  "

Set-TextValue 10 4 "59.71"
Set-TextValue 10 5 "  +0.21%  "

Set-TextValue 11 4 "0.0780"
Set-TextValue 11 5 "  +5.13%  "

Set-TextValue 12 5 "  +0.08%  "

Set-TextValue 13 5 "  -2.75%  "

Set-TextValue 14 4 "23.76"
Set-TextValue 14 5 "  +26.16%  "

Set-TextValue 15 4 "14.83"
Set-TextValue 15 5 "  +0.72%  "

Set-TextValue 16 4 "2.346.11"
Set-TextValue 16 5 "  +3.05%  "

Set-TextValue 17 4 "5.71"
Set-TextValue 17 5 "  +6.99%  "

Set-TextValue 18 4 "2.041.04"
Set-TextValue 18 5 "  +2.69%  "

Set-TextValue 19 4 "37.223.24"
Set-TextValue 19 5 "  +4.19%  "

Set-TextValue 20 4 "73.60"
Set-TextValue 20 5 "  +2.47%  "

Set-TextValue 21 4 "0.0₃0885"
Set-TextValue 21 5 "  +3.78%  "

Set-TextValue 22 4 "5.49"
Set-TextValue 22 5 "  +5.29%  "

Set-TextValue 23 4 "239.86"
Set-TextValue 23 5 "  +2.70%  "

Set-TextValue 24 2 "PancakeSwap"
Set-TextValue 24 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue 24 4 "2.63"
Set-TextValue 24 5 "  +1.81%  "

Set-TextValue 25 2 "Dai"
Set-TextValue 25 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue 25 4 "1.00"
Set-TextValue 25 5 "  -0.13%  "

Set-TextValue 26 5 "  +4.79%  "

Set-TextValue 27 4 "10.21"
Set-TextValue 27 5 "  +9.79%  "

Set-TextValue 28 4 "161.83"
Set-TextValue 28 5 "  -1.77%  "

Set-TextValue 29 2 "Kaspa"
Set-TextValue 29 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue 29 4 "0.135"
Set-TextValue 29 5 "  +37.17%  "

Set-TextValue 30 2 "EthereumClassic"
Set-TextValue 30 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue 30 4 "20.07"
Set-TextValue 30 5 "  +3.84%  "

Set-TextValue 31 5 "  +2.64%  "

Set-TextValue 32 4 "5.19"
Set-TextValue 32 5 "  +5.42%  "

Set-TextValue 33 5 "  +6.01%  "

Set-TextValue 34 2 "Hedera"
Set-TextValue 34 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue 34 4 "0.0630"
Set-TextValue 34 5 "  +4.70%  "

Set-TextValue 35 2 "InternetComputer(DFINITY)"
Set-TextValue 35 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue 35 4 "4.70"
Set-TextValue 35 5 "  +7.03%  "

Set-TextValue 36 4 "2.40"
Set-TextValue 36 5 "  -2.72%  "

Set-TextValue 37 4 "6.36"
Set-TextValue 37 5 "  +11.64%  "

Set-TextValue 38 5 "  -0.13%  "

Set-TextValue 39 5 "  +2.64%  "

Set-TextValue 40 4 "3.06"
Set-TextValue 40 5 "  +31.40%  "

Set-TextValue 41 4 "0.102"
Set-TextValue 41 5 "  +8.09%  "

Set-TextValue 42 4 "1.28"
Set-TextValue 42 5 "  +3.92%  "

Set-TextValue 43 4 "3.06"
Set-TextValue 43 5 "  +7.28%  "

Set-TextValue 44 4 "17.71"
Set-TextValue 44 5 "  +7.10%  "

Set-TextValue 45 5 "  +5.78%  "

Set-TextValue 46 4 "0.0220"
Set-TextValue 46 5 "  +2.63%  "

Set-TextValue 47 4 "96.42"
Set-TextValue 47 5 "  +2.45%  "

Set-TextValue 48 4 "7.85"
Set-TextValue 48 5 "  +1.00%  "

Set-TextValue 49 4 "1.405.30"
Set-TextValue 49 5 "  +2.63%  "

Set-TextValue 50 4 "2.94"
Set-TextValue 50 5 "  +1.51%  "

Set-TextValue 51 4 "47.72"
Set-TextValue 51 5 "  +1.14%  "

